# Weekly fruit/vegetable update: two new price records were logged for
# Terminal Hortofrutícola Agro Chillán - Choclo (Provincia de Diguillín,
# 2023-04-05), pushing the previously-existing rows 290-309 down to 292-311.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows right above the old row 290 (each Insert() pushes
# the row that's currently at 290, plus everything below it, down by one -
# doing it twice opens up two blank rows at 290 and 291).
$ws.Rows.Item(290).Insert()
$ws.Rows.Item(290).Insert()

# Common columns shared by every data row in this sheet.
$mercadoId = 7
$mercado = "Terminal Hortofrutícola Agro Chillán"
$region = "Ñuble"
$codreg = 16
$categoriaId = 100112024
$categoria = "Choclo"
$clasificacion = "Hortaliza"

# New row 290
$ws.Range("A290").Value = $mercadoId
$ws.Range("B290").Value = $mercado
$ws.Range("C290").Value = $region
$ws.Range("D290").Value = 45021
$ws.Range("E290").Value = $codreg
$ws.Range("F290").Value = $categoriaId
$ws.Range("G290").Value = $categoria
$ws.Range("H290").Value = "Choclero"
$ws.Range("I290").Value = "Primera"
$ws.Range("J290").Value = 1800
$ws.Range("K290").Value = 450
$ws.Range("L290").Value = 500
$ws.Range("M290").Value = 472
$ws.Range("N290").Value = "$/unidad"
$ws.Range("O290").Value = "Provincia de Diguillín"
$ws.Range("P290").Value = 472
$ws.Range("Q290").Value = 1
$ws.Range("R290").Value = $clasificacion

# New row 291
$ws.Range("A291").Value = $mercadoId
$ws.Range("B291").Value = $mercado
$ws.Range("C291").Value = $region
$ws.Range("D291").Value = 45021
$ws.Range("E291").Value = $codreg
$ws.Range("F291").Value = $categoriaId
$ws.Range("G291").Value = $categoria
$ws.Range("H291").Value = "Choclero"
$ws.Range("I291").Value = "Segunda"
$ws.Range("J291").Value = 400
$ws.Range("K291").Value = 350
$ws.Range("L291").Value = 350
$ws.Range("M291").Value = 350
$ws.Range("N291").Value = "$/unidad"
$ws.Range("O291").Value = "Provincia de Diguillín"
$ws.Range("P291").Value = 350
$ws.Range("Q291").Value = 1
$ws.Range("R291").Value = $clasificacion
